# "add download excel button funtionality"
# Target edit (per xml diff):
#  - add a new worksheet "bla" after Sheet1 with a small lookup table
#    (some list / number) used for data entry / validation
#  - highlight row 2 of Sheet1 with a light "Accent1 Lighter 80%" fill
#  - add a handful of blank formatted rows below it (4,6,8,10) to give
#    room for more rows of entered data
#  - add a whole-number (0-100) data validation on Sheet1!A2:A10

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Add the "bla" lookup sheet right after Sheet1
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "bla"

$colA = @("some list", "asdf", "asdf", "sadf", "asdf", "asdf", "sda", "f")
$colB = @("number", 1, 2, 3, 4, 4, 5, 5)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws2.Cells.Item($i + 1, 2).Value = $colB[$i]
}

$ws2.Range("B7").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Go back to Sheet1 and highlight row 2 (light blue fill -
#    "Blue, Accent 1, Lighter 80%")
# ---------------------------------------------------------------------
$ws1.Activate()

$highlightColor = 16247774   # RGB(0xDE,0xEB,0xF7)
$row2Cols = @("A2","B2","C2","D2","E2","F2","G2","I2","K2","L2","M2","N2","O2", `
              "V2","W2","X2","Y2","Z2","AB2","AC2","AD2","AF2","AG2","AH2","AI2", `
              "AJ2","AK2","AL2","AM2","AO2","AP2","AQ2","AU2","AY2","AZ2","BA2","BB2")
foreach ($addr in $row2Cols) {
    $ws1.Range($addr).Interior.Color = $highlightColor
}

# A few more (currently empty) rows get the same row formatting so the
# user has extra pre-formatted rows ready for new data
foreach ($r in 4, 6, 8, 10) {
    $ws1.Cells.Item($r, 1).Interior.Color = $highlightColor
}

# ---------------------------------------------------------------------
# 3) Whole-number data validation (0-100) on column A, rows 2-10
# ---------------------------------------------------------------------
$valRange = $ws1.Range("A2:A10")
$valRange.Validation.Delete()
$valRange.Validation.Add(1, 1, 1, 0, 100)
$valRange.Validation.IgnoreBlank = $true
$valRange.Validation.InCellDropdown = $true
$valRange.Validation.ShowInput = $true
$valRange.Validation.ShowError = $true

$ws1.Range("A1").Select() | Out-Null

Write-Output "done"
